$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.741185426712036
$ws.Range("B1").Value = 2.238624811172485
$ws.Range("C1").Value = 2.393944501876831
$ws.Range("D1").Value = 3.130728483200073
$ws.Range("E1").Value = 1.81594979763031
